$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A1:A3 values
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3

# Clear C1 (was a shared string, now removed) and B3 (formula removed)
$ws.Range("C1").Clear()
$ws.Range("B3").Clear()

# B1 becomes a formula
$ws.Range("B1").Formula = "=A1*2"
